$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 140.2
$ws.Range("I5").Value = 140.2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 140.2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -25.19999999999999
$ws.Range("N5").ClearContents()
# Row 98
$ws.Range("H98").Value = 2030.7307
$ws.Range("I98").Value = 2022.6818
$ws.Range("K98").Value = 2022.6818
$ws.Range("M98").Value = -524.6818000000001
# Row 122
$ws.Range("H122").Value = 2030.7307
$ws.Range("I122").Value = 2022.6818
$ws.Range("K122").Value = 6068.0454
$ws.Range("M122").Value = -3618.0454
# Row 132
$ws.Range("H132").Value = 2441368
$ws.Range("I132").Value = 2271.375
$ws.Range("J132").Value = 11113711
$ws.Range("K132").Value = 6814.125
$ws.Range("L132").Value = 33341133
$ws.Range("M132").Value = -4284.125
$ws.Range("N132").Value = -33346193
# Row 133
$ws.Range("H133").Value = 48523.08
$ws.Range("I133").Value = 25800
$ws.Range("J133").Value = 52654.547
$ws.Range("K133").Value = 25800
$ws.Range("L133").Value = 52654.547
$ws.Range("M133").Value = -20740
$ws.Range("N133").Value = -62774.547
# Row 137
$ws.Range("H137").Value = 1056.75
$ws.Range("I137").Value = 838.2143
$ws.Range("J137").Value = 1566.6666
$ws.Range("K137").Value = 2514.6429
$ws.Range("L137").Value = 4699.9998
$ws.Range("M137").Value = 35.35710000000017
$ws.Range("N137").Value = -9799.9998
# Row 138
$ws.Range("H138").Value = 3967
$ws.Range("I138").Value = 3118.3076
$ws.Range("J138").Value = 4095.2908
$ws.Range("K138").Value = 9354.9228
$ws.Range("L138").Value = 12285.8724
$ws.Range("M138").Value = -4214.9228
$ws.Range("N138").Value = -22565.8724

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17768.773
$ws.Range("I32").Value = 12776.122
$ws.Range("J32").Value = 54986.727
$ws.Range("K32").Value = 12776.122
$ws.Range("L32").Value = 54986.727
$ws.Range("M32").Value = -12489.122
$ws.Range("N32").Value = -55560.727
# Row 74
$ws.Range("H74").Value = 1168.0513
$ws.Range("I74").Value = 1166.3226
$ws.Range("J74").Value = 1174.75
$ws.Range("K74").Value = 1166.3226
$ws.Range("L74").Value = 1174.75
$ws.Range("M74").Value = -292.3226
$ws.Range("N74").Value = -2922.75
# Row 77
$ws.Range("H77").Value = 1168.0513
$ws.Range("I77").Value = 1166.3226
$ws.Range("J77").Value = 1174.75
$ws.Range("K77").Value = 5831.612999999999
$ws.Range("L77").Value = 5873.75
$ws.Range("M77").Value = -1463.612999999999
$ws.Range("N77").Value = -14609.75
# Row 132
$ws.Range("H132").Value = 2468.3333
$ws.Range("I132").Value = 1852.15
$ws.Range("J132").Value = 3700.7
$ws.Range("K132").Value = 5556.450000000001
$ws.Range("L132").Value = 11102.1
$ws.Range("M132").Value = -3026.450000000001
$ws.Range("N132").Value = -16162.1

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 15018.797
$ws.Range("I134").Value = 1297.5245
$ws.Range("J134").Value = 79403.234
$ws.Range("K134").Value = 3892.5735
$ws.Range("L134").Value = 238209.702
$ws.Range("M134").Value = -1357.5735
$ws.Range("N134").Value = -243279.702

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2537.761
$ws.Range("I31").Value = 2126
$ws.Range("J31").Value = 4020.1
$ws.Range("K31").Value = 2126
$ws.Range("L31").Value = 4020.1
$ws.Range("M31").Value = -1831
$ws.Range("N31").Value = -4610.1
# Row 34
$ws.Range("H34").Value = 2537.761
$ws.Range("I34").Value = 2126
$ws.Range("J34").Value = 4020.1
$ws.Range("K34").Value = 2126
$ws.Range("L34").Value = 4020.1
$ws.Range("M34").Value = -1924
$ws.Range("N34").Value = -4424.1
# Row 132
$ws.Range("H132").Value = 1506.9678
$ws.Range("I132").Value = 1139.1765
$ws.Range("J132").Value = 1953.5714
$ws.Range("K132").Value = 3417.5295
$ws.Range("L132").Value = 5860.7142
$ws.Range("M132").Value = -887.5295000000001
$ws.Range("N132").Value = -10920.7142
# Row 134
$ws.Range("H134").Value = 50001444
$ws.Range("I134").Value = 1730.7142
$ws.Range("J134").Value = 166667440
$ws.Range("K134").Value = 5192.142599999999
$ws.Range("L134").Value = 500002320
$ws.Range("M134").Value = -2657.142599999999
$ws.Range("N134").Value = -500007390

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 124.75
$ws.Range("J23").Value = 145.54546
$ws.Range("L23").Value = 436.63638
$ws.Range("N23").Value = -906.6363799999999
# Row 116
$ws.Range("H116").Value = 2422.2222
$ws.Range("I116").Value = 2280
$ws.Range("J116").Value = 2600
$ws.Range("K116").Value = 6840
$ws.Range("L116").Value = 7800
$ws.Range("M116").Value = -3398
$ws.Range("N116").Value = -14684
# Row 118
$ws.Range("H118").Value = 5414.5
$ws.Range("I118").Value = 4029
$ws.Range("J118").Value = 6800
$ws.Range("K118").Value = 12087
$ws.Range("L118").Value = 20400
$ws.Range("M118").Value = -10844
$ws.Range("N118").Value = -22886
# Row 129
$ws.Range("H129").Value = 29562.184
$ws.Range("I129").Value = 1294.9
$ws.Range("J129").Value = 39657.645
$ws.Range("K129").Value = 3884.7
$ws.Range("L129").Value = 118972.935
$ws.Range("M129").Value = 1115.3
$ws.Range("N129").Value = -128972.935
# Row 131
$ws.Range("H131").Value = 78522.08
$ws.Range("I131").Value = 101457
$ws.Range("K131").Value = 304371
$ws.Range("M131").Value = -299331
# Row 137
$ws.Range("H137").Value = 86553.336
$ws.Range("I137").Value = 3314
$ws.Range("K137").Value = 9942
$ws.Range("M137").Value = -4842

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 7218.6665
$ws.Range("I126").Value = 4662.4
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 13987.2
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -11517.2
$ws.Range("N126").Value = -64940

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 811.0909
$ws.Range("J22").Value = 833.6667
$ws.Range("L22").Value = 833.6667
$ws.Range("N22").Value = -1423.6667
# Row 27
$ws.Range("H27").Value = 811.0909
$ws.Range("J27").Value = 833.6667
$ws.Range("L27").Value = 833.6667
$ws.Range("N27").Value = -1047.6667
# Row 141
$ws.Range("H141").Value = 72700
$ws.Range("J141").Value = 72700
$ws.Range("L141").Value = 72700
$ws.Range("N141").Value = -83060

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 9800
$ws.Range("I54").Value = 9800
$ws.Range("K54").Value = 9800
$ws.Range("M54").Value = -9280
